# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    3 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    4 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    5 = @(1.445647641019636, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 2.741163836583335)
    6 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
